$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "quantity" column (F) values for the affected rows.
$ws.Range("F2").Value = -239
$ws.Range("F3").Value = -530
$ws.Range("F4").Value = -832
$ws.Range("F5").Value = -107
$ws.Range("F7").Value = -107
